$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the header labels: "<Name>_old" -> "<Name>_FV2310", "<Name>_new" -> "<Name>_FV2404"
$headers = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = "$($headers[$i])_FV2310"
}

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = "$($headers[$i])_FV2404"
}

# 2. Freeze the header row (split below row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Turn the used range into an Excel Table ("Table1").
$rng = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
